# Add a new acronym-key row for "AVIC" (Annual Vehicle Insurance Cost)
# to the "Key to Variables" sheet, just above the existing "AVL" row,
# and make "Key to Variables" the active sheet tab (it was "About" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# Insert a fresh row at 187 (pushes the old row 187 "AVL" and everything
# below it down by one row).
$ws.Rows.Item(187).Insert()

# Row 168 ("AVLo" / Average Vehicle Loading) already has the formatting
# we want for this new row: plain wrap-text cells for A:E & G, and the
# pale-yellow "medium importance" fill (style used on column F) - so
# copy its formatting into the new row before filling in values.
$ws.Range("A168:G168").Copy()
$ws.Range("A187:G187").PasteSpecial(-4122)

# Fill in the new row's content.
$ws.Cells.Item(187, 1).Value = "trans"
$ws.Cells.Item(187, 2).Value = "AVIC"
$ws.Cells.Item(187, 3).Value = "Annual Vehicle Insurance Cost"
$ws.Cells.Item(187, 6).Value = "medium"

# Make "Key to Variables" the active/selected worksheet tab.
$ws.Activate()
